# Applies the "Updated cryptos list" price/volume refresh described by the commit.
# Source values are plain numeric/text strings pulled from coinranking.com and stored
# as literal text in columns D (Price) and E (Volume/1h) -- never as real numbers -- so
# values like "141.00" or "1.009" keep their exact printed digits (incl. trailing zeros)
# instead of being normalised by Excel's numeric parser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold strings that Excel would otherwise auto-detect as numbers
# (e.g. "1.009", "141.00") and silently reformat/trim -- pre-mark them as Text so the
# literal string is preserved verbatim, matching how the rest of the sheet stores them.
$textCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D18", "D19",
    "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34",
    "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48",
    "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Write every changed cell (Coin / Link / Price / Volume(1h)) row by row.

$ws.Range("D2").Value = "20.513.06"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").Value = "1.473.68"
$ws.Range("E3").Value = "  +3.49%  "

$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.84%  "

$ws.Range("D5").Value = "0.9575"
$ws.Range("E5").Value = "  -3.95%  "

$ws.Range("D6").Value = "276.82"
$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").Value = "0.3654"
$ws.Range("E7").Value = "  -1.52%  "

$ws.Range("D8").Value = "0.3060"
$ws.Range("E8").Value = "  -2.76%  "

$ws.Range("D9").Value = "39.68"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Value = "1.054"
$ws.Range("E10").Value = "  -1.10%  "

$ws.Range("D11").Value = "0.06609"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").Value = "18.11"
$ws.Range("E13").Value = "  +0.58%  "

$ws.Range("D14").Value = "5.456"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").Value = "6.178"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").Value = "1.474.64"
$ws.Range("E17").Value = "  +3.64%  "

$ws.Range("D18").Value = "0.05901"
$ws.Range("E18").Value = "  +3.12%  "

$ws.Range("D19").Value = "0.9632"
$ws.Range("E19").Value = "  -3.42%  "

$ws.Range("D20").Value = "69.07"
$ws.Range("E20").Value = "  -3.94%  "

$ws.Range("D21").Value = "5.467"
$ws.Range("E21").Value = "  -2.75%  "

$ws.Range("D22").Value = "14.47"
$ws.Range("E22").Value = "  -2.90%  "

$ws.Range("D23").Value = "11.03"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("D24").Value = "2.260"
$ws.Range("E24").Value = "  +1.56%  "

$ws.Range("D25").Value = "20.566.65"
$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("D26").Value = "141.00"
$ws.Range("E26").Value = "  +4.67%  "

$ws.Range("D27").Value = "2.119"
$ws.Range("E27").Value = "  -7.93%  "

$ws.Range("D28").Value = "17.17"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").Value = "1.631.79"
$ws.Range("E29").Value = "  +3.15%  "

$ws.Range("D30").Value = "113.53"
$ws.Range("E30").Value = "  +2.08%  "

$ws.Range("D31").Value = "3.949"
$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").Value = "0.8181"
$ws.Range("E32").Value = "  -1.55%  "

$ws.Range("D33").Value = "4.959"
$ws.Range("E33").Value = "  -6.47%  "

$ws.Range("D34").Value = "0.07943"
$ws.Range("E34").Value = "  +1.58%  "

$ws.Range("E35").Value = "  +3.52%  "

$ws.Range("D36").Value = "1.225"
$ws.Range("E36").Value = "  +10.29%  "

$ws.Range("D37").Value = "0.05766"
$ws.Range("E37").Value = "  -1.85%  "

$ws.Range("D38").Value = "4.728"
$ws.Range("E38").Value = "  -4.04%  "

$ws.Range("D39").Value = "0.02035"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").Value = "10.42"
$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "7.587"
$ws.Range("E41").Value = "  -4.99%  "

$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "0.9561"
$ws.Range("E42").Value = "  -4.02%  "

$ws.Range("D43").Value = "0.1877"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").Value = "0.5274"
$ws.Range("E44").Value = "  -1.63%  "

$ws.Range("D45").Value = "3.503"
$ws.Range("E45").Value = "  -1.44%  "

$ws.Range("D46").Value = "12.00"
$ws.Range("E46").Value = "  -2.69%  "

$ws.Range("D47").Value = "117.40"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("D48").Value = "0.5187"
$ws.Range("E48").Value = "  -1.25%  "

$ws.Range("D49").Value = "1.777"
$ws.Range("E49").Value = "  -0.63%  "

$ws.Range("D50").Value = "0.06455"

$ws.Range("D51").Value = "0.9960"
$ws.Range("E51").Value = "  -0.11%  "
